# 🚌 141: 30/12 14:47 LP1912+6203+6173
# Appends freshly scraped rows to the three schedule sheets and refreshes
# the "last updated" / "total rows" header cells on each sheet.

$wb = $excel.ActiveWorkbook

$timestamp = "30/12/2025 11:47:44"

# ---------------------------------------------------------------------
# Sheet "LP1912": columns B..G = Hora_Scrap, Hora_Llegada, Linea, Minutos,
# Parada, Fecha  (column A is always blank on data rows)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $timestamp"
$ws1.Range("A3").Value = "Total filas: 187"

$sheet1Rows = @(
    @("11:47:33","11:51","215B_EL PATO",4,"LP1912","30/12/2025"),
    @("11:47:33","11:52","15_ABASTO",5,"LP1912","30/12/2025"),
    @("11:47:33","11:56","16_SANTA ANA",9,"LP1912","30/12/2025"),
    @("11:47:33","12:02","84_COLONIA URQUIZA-ESC 49",15,"LP1912","30/12/2025"),
    @("11:47:33","12:06","16_P MOR-SANTA ANA",19,"LP1912","30/12/2025"),
    @("11:47:33","12:07","23_HERNANDEZ",20,"LP1912","30/12/2025"),
    @("11:47:33","12:13","10_OLMOS",26,"LP1912","30/12/2025"),
    @("11:47:33","12:16","16_SANTA ANA",29,"LP1912","30/12/2025"),
    @("11:47:33","12:20","14_ABASTO",33,"LP1912","30/12/2025"),
    @("11:47:33","12:21","26_HERNANDEZ",34,"LP1912","30/12/2025"),
    @("11:47:33","12:34","23_HERNANDEZ",47,"LP1912","30/12/2025"),
    @("11:47:33","12:38","17_179 Y 38",51,"LP1912","30/12/2025"),
    @("11:47:33","12:48","11_ETCHEVERRY",61,"LP1912","30/12/2025"),
    @("11:47:33","12:50","15_ABASTO",63,"LP1912","30/12/2025"),
    @("11:47:33","12:55","10_OLMOS",68,"LP1912","30/12/2025"),
    @("11:47:33","13:06","16_P MOR-SANTA ANA",79,"LP1912","30/12/2025"),
    @("11:47:33","13:16","17_ROMERO",89,"LP1912","30/12/2025")
)

$row = 172
foreach ($r in $sheet1Rows) {
    $ws1.Cells.Item($row, 2).Value = $r[0]
    $ws1.Cells.Item($row, 3).Value = $r[1]
    $ws1.Cells.Item($row, 4).Value = $r[2]
    $ws1.Cells.Item($row, 5).Value = $r[3]
    $ws1.Cells.Item($row, 6).Value = $r[4]
    $ws1.Cells.Item($row, 7).Value = $r[5]
    $row++
}

# ---------------------------------------------------------------------
# Sheet "LP1912-215": columns B..G = Fecha, Hora_Scrap, Hora_Llegada,
# Linea, Minutos, Parada
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $timestamp"
$ws2.Range("A3").Value = "Total filas: 18"

$ws2.Cells.Item(19, 2).Value = "30/12/2025"
$ws2.Cells.Item(19, 3).Value = "11:47:33"
$ws2.Cells.Item(19, 4).Value = "11:51"
$ws2.Cells.Item(19, 5).Value = "215B_EL PATO"
$ws2.Cells.Item(19, 6).Value = 4
$ws2.Cells.Item(19, 7).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet "6203-6173": columns B..G = Fecha, Hora_Scrap, Hora_Llegada,
# Linea, Minutos, Parada
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $timestamp"
$ws3.Range("A3").Value = "Total filas: 25"

$ws3.Cells.Item(25, 2).Value = "30/12/2025"
$ws3.Cells.Item(25, 3).Value = "11:47:44"
$ws3.Cells.Item(25, 4).Value = "12:04"
$ws3.Cells.Item(25, 5).Value = "215A_LA PLATA"
$ws3.Cells.Item(25, 6).Value = 17
$ws3.Cells.Item(25, 7).Value = "L6173"

$ws3.Cells.Item(26, 2).Value = "30/12/2025"
$ws3.Cells.Item(26, 3).Value = "11:47:39"
$ws3.Cells.Item(26, 4).Value = "12:53"
$ws3.Cells.Item(26, 5).Value = "215C_LA PLATA"
$ws3.Cells.Item(26, 6).Value = 66
$ws3.Cells.Item(26, 7).Value = "L6203"
